$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (price + 1h volume%), plus two rank swaps:
#   row17<->row18 (Chainlink <-> TRON) and row19<->row20 (Uniswap <-> WrappedBTC).
#
# The "Price" column holds plain text (not real numbers: values like
# "68.265.54" use dots as thousands separators and would misparse as a
# number/date). Force text storage for every Price write: flip the cell to
# the "@" (Text) format before assigning .Value, then restore the "Normal"
# style so no stray number-format/quote-prefix is left behind on the cell.
function Set-PriceText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-PriceText "D17" '0.127'
$ws.Range("E17").Value = '  +1.69%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-PriceText "D18" '19.29'
$ws.Range("E18").Value = '  +10.24%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-PriceText "D19" '68.131.74'
$ws.Range("E19").Value = '  +5.38%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-PriceText "D20" '12.42'
$ws.Range("E20").Value = '  +8.32%  '

Set-PriceText "D2" '68.265.54'
$ws.Range("E2").Value = '  +5.12%  '
Set-PriceText "D3" '3.618.41'
$ws.Range("E3").Value = '  +5.35%  '
$ws.Range("E4").Value = '  +0.09%  '
Set-PriceText "D5" '200.89'
$ws.Range("E5").Value = '  +11.69%  '
Set-PriceText "D6" '577.46'
$ws.Range("E6").Value = '  +4.41%  '
Set-PriceText "D7" '3.612.73'
$ws.Range("E7").Value = '  +5.39%  '
Set-PriceText "D8" '0.621'
$ws.Range("E8").Value = '  +5.01%  '
$ws.Range("E9").Value = '  -0.34%  '
Set-PriceText "D10" '0.683'
$ws.Range("E10").Value = '  +7.81%  '
Set-PriceText "D11" '60.44'
Set-PriceText "D12" '0.149'
$ws.Range("E12").Value = '  +7.94%  '
Set-PriceText "D13" '0.0000284'
$ws.Range("E13").Value = '  +15.71%  '
Set-PriceText "D14" '10.22'
$ws.Range("E14").Value = '  +10.21%  '
Set-PriceText "D15" '4.202.20'
$ws.Range("E15").Value = '  +5.55%  '
Set-PriceText "D16" '3.627.61'
$ws.Range("E16").Value = '  +5.40%  '
Set-PriceText "D21" '1.07'
$ws.Range("E21").Value = '  +6.17%  '
Set-PriceText "D22" '405.87'
$ws.Range("E22").Value = '  +9.01%  '
Set-PriceText "D23" '12.91'
$ws.Range("E23").Value = '  +23.06%  '
Set-PriceText "D24" '4.23'
$ws.Range("E24").Value = '  +5.40%  '
Set-PriceText "D25" '85.58'
$ws.Range("E25").Value = '  +5.02%  '
Set-PriceText "D26" '4.00'
$ws.Range("E26").Value = '  +19.35%  '
Set-PriceText "D27" '2.93'
$ws.Range("E27").Value = '  +6.80%  '
Set-PriceText "D28" '12.61'
$ws.Range("E28").Value = '  +8.36%  '
Set-PriceText "D29" '6.14'
$ws.Range("E29").Value = '  +2.75%  '
Set-PriceText "D30" '9.36'
$ws.Range("E30").Value = '  +11.98%  '
Set-PriceText "D31" '7.79'
$ws.Range("E31").Value = '  +12.43%  '
Set-PriceText "D32" '31.82'
$ws.Range("E32").Value = '  +7.09%  '
Set-PriceText "D33" '680.98'
$ws.Range("E33").Value = '  +14.01%  '
Set-PriceText "D34" '12.23'
$ws.Range("E34").Value = '  +5.56%  '
Set-PriceText "D35" '0.114'
$ws.Range("E35").Value = '  +6.12%  '
Set-PriceText "D36" '63.77'
$ws.Range("E36").Value = '  +2.53%  '
Set-PriceText "D37" '41.78'
$ws.Range("E37").Value = '  +5.58%  '
Set-PriceText "D38" '0.415'
$ws.Range("E38").Value = '  +8.37%  '
$ws.Range("E39").Value = '  +0.13%  '
Set-PriceText "D40" '0.0₃0768'
$ws.Range("E40").Value = '  +9.87%  '
Set-PriceText "D41" '3.19'
$ws.Range("E41").Value = '  +19.68%  '
Set-PriceText "D42" '3.198.70'
$ws.Range("E42").Value = '  +10.67%  '
Set-PriceText "D43" '0.135'
$ws.Range("E43").Value = '  +7.46%  '
Set-PriceText "D44" '0.999'
$ws.Range("E44").Value = '  +0.11%  '
Set-PriceText "D45" '2.69'
$ws.Range("E45").Value = '  +13.05%  '
Set-PriceText "D46" '2.86'
$ws.Range("E46").Value = '  +30.13%  '
Set-PriceText "D47" '2.80'
$ws.Range("E47").Value = '  +14.91%  '
Set-PriceText "D48" '0.0417'
$ws.Range("E48").Value = '  +8.13%  '
$ws.Range("E49").Value = '  +5.48%  '
Set-PriceText "D50" '8.77'
$ws.Range("E50").Value = '  +10.09%  '
$ws.Range("E51").Value = '  +0.96%  '
